$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# --- Row 27: new time-registration entry ---
$ws.Range("A27").Value = "DD + DD07 sammenlægning og  samlet dataordbog oprettet"
$ws.Range("B27").Value = "business-Process Analyst"
$ws.Range("C27").Value = 43896
$ws.Range("D27").Value = 0.36458333333333331
$ws.Range("E27").Value = 0.45833333333333331

# --- Row 28: new time-registration entry ---
$ws.Range("A28").Value = "ATD03a + ATD03b Beregn Bruttofortjeneste sammenlægning"
$ws.Range("B28").Value = "business-Process Analyst"
$ws.Range("C28").Value = 43896
$ws.Range("D28").Value = 0.46180555555555558
$ws.Range("E28").Value = 0.54652777777777783

# --- Row 29: new time-registration entry ---
$ws.Range("A29").Value = "ATD06 KKO sammenlægning"
$ws.Range("B29").Value = "business-Process Analyst"
$ws.Range("C29").Value = 43896
$ws.Range("D29").Value = 0.54791666666666672
$ws.Range("E29").Value = 0.56736111111111109

# --- Row 30: new time-registration entry ---
$ws.Range("A30").Value = "Opfølgning og tilføjelser af sammenlægning af DDér i Ordbogen"
$ws.Range("B30").Value = "business-Process Analyst"
$ws.Range("C30").Value = 43896
$ws.Range("D30").Value = 0.57291666666666663
$ws.Range("E30").Value = 0.63541666666666663

# Move the active selection to B30, matching where the author ended up editing
$ws.Range("B30").Select()
